$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3770616352558136
$ws.Range("B1").Value = 1.102589726448059
$ws.Range("C1").Value = 4.730323314666748
$ws.Range("D1").Value = 1.799968600273132
$ws.Range("E1").Value = 1.002940535545349
